$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("D2").Value = 44420
$ws.Range("H2").Value = 'Madrigal'
$ws.Range("J2").Value = 800
$ws.Range("K2").Value = 14000
$ws.Range("L2").Value = 15000
$ws.Range("M2").Value = 14500
$ws.Range("N2").Value = '$/caja 40 unidades'
$ws.Range("O2").Value = 'Provincia de Limarí'
$ws.Range("P2").Value = 362
$ws.Range("Q2").Value = 40

# Row 3
$ws.Range("J3").Value = 700
$ws.Range("K3").Value = 13000
$ws.Range("L3").Value = 14000
$ws.Range("M3").Value = 13500
$ws.Range("O3").Value = 'Provincia del Elquí'
$ws.Range("P3").Value = 338

# Row 4
$ws.Range("D4").Value = 44427
$ws.Range("J4").Value = 400
$ws.Range("K4").Value = 12000
$ws.Range("L4").Value = 13000
$ws.Range("M4").Value = 12500
$ws.Range("O4").Value = 'Provincia de Limarí'
$ws.Range("P4").Value = 312

# Row 5
$ws.Range("D5").Value = 44498
$ws.Range("K5").Value = 8500
$ws.Range("L5").Value = 9000
$ws.Range("M5").Value = 8750
$ws.Range("O5").Value = 'Provincia de Limarí'
$ws.Range("P5").Value = 292

# Row 6
$ws.Range("D6").Value = 44484
$ws.Range("J6").Value = 300
$ws.Range("K6").Value = 9000
$ws.Range("L6").Value = 10000
$ws.Range("M6").Value = 9500
$ws.Range("O6").Value = 'Provincia del Elquí'
$ws.Range("P6").Value = 317

# Row 9
$ws.Range("D9").Value = 44438
$ws.Range("H9").Value = 'Española'
$ws.Range("K9").Value = 11000
$ws.Range("L9").Value = 12000
$ws.Range("M9").Value = 11500
$ws.Range("N9").Value = '$/caja 30 unidades'
$ws.Range("O9").Value = 'Provincia del Elquí'
$ws.Range("P9").Value = 383
$ws.Range("Q9").Value = 30
